{"js": "// The \"Avalia\u00e7\u00e3o\" section currently has two runs whose text concatenates\n// several sentences without any separator between them. The edit breaks\n// each of those runs apart at the sentence boundaries by inserting a\n// manual line break (<w:br/>, i.e. the \"\\u000b\" vertical-tab character in\n// the Word text model) right after each sentence.\n//\n// Strategy: for every boundary, search() for the (unique) text that ends\n// exactly at the boundary, then insert a line-break character immediately\n// \"After\" that found range. Office.js's Range.insertText(\"\\u000b\", \"After\")\n// lands a real <w:br/> at that exact point, splitting the run's <w:t> into\n// two <w:t> elements around the new <w:br/> - matching how Word itself\n// represents a manual line break.\n\nasync function insertLineBreakAfter(anchorText) {\n  const results = context.document.body.search(anchorText, { matchCase: true });\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find anchor text: \" + anchorText);\n  }\n\n  results.items[0].insertText(\"\\u000b\", \"After\");\n  await context.sync();\n}\n\n// \"Crit\u00e9rio:\" run -> split into the two sentences it concatenates.\nawait insertLineBreakAfter(\"Para os estudantes: despertar interesse na engenharia.\");\n\n// \"Norma de recupera\u00e7\u00e3o:\" run -> split into its five \"- ...\" bullet sentences.\nawait insertLineBreakAfter(\"necessidades, desafios e prefer\u00eancias dos estudantes.\");\nawait insertLineBreakAfter(\"e quaisquer outras considera\u00e7\u00f5es importantes.\");\nawait insertLineBreakAfter(\"sobre a profiss\u00e3o engenharia.\");\nawait insertLineBreakAfter(\"conhecimento sobre o tema.\");\n", "ps1": "# The \"Avalia\u00e7\u00e3o\" section currently has two runs whose text concatenates\n# several sentences without any separator between them. The edit breaks\n# each of those runs apart at the sentence boundaries by inserting a\n# manual line break (a Word \"^l\" / vertical-tab break, i.e. <w:br/>) right\n# after each sentence.\n#\n# Strategy: use Find/Replace (wdFindContinue, no wildcards needed) on the\n# whole document, replacing \"<end of sentence><start of next sentence>\"\n# with \"<end of sentence>^l<start of next sentence>\" - ^l is Word's\n# replacement-text code for a manual line break. Each find string is a\n# short, unique snippet spanning exactly the boundary being split.\n\n$d = $word.ActiveDocument\n\nfunction Insert-LineBreak($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $ok = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $findText\"\n    }\n}\n\n# \"Crit\u00e9rio:\" run -> split into the two sentences it concatenates.\nInsert-LineBreak \"despertar interesse na engenharia.Para a forma\u00e7\u00e3o\" \"despertar interesse na engenharia.^lPara a forma\u00e7\u00e3o\"\n\n# \"Norma de recupera\u00e7\u00e3o:\" run -> split into its five \"- ...\" bullet sentences.\nInsert-LineBreak \"dos estudantes.- Defini\u00e7\u00e3o\" \"dos estudantes.^l- Defini\u00e7\u00e3o\"\nInsert-LineBreak \"outras considera\u00e7\u00f5es importantes.- Pesquisa\" \"outras considera\u00e7\u00f5es importantes.^l- Pesquisa\"\nInsert-LineBreak \"sobre a profiss\u00e3o engenharia.- Avalia\u00e7\u00e3o\" \"sobre a profiss\u00e3o engenharia.^l- Avalia\u00e7\u00e3o\"\nInsert-LineBreak \"conhecimento sobre o tema.- Implementa\u00e7\u00e3o\" \"conhecimento sobre o tema.^l- Implementa\u00e7\u00e3o\"\n"}
